# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2-210) from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C210").Value = 45175

Write-Host "Updated C2:C210 to 45175"
